$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "45.856.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.94%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.375.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.15%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.87%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.564"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.36%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  -3.75%  "
$ws.Range("E10").Value = "  -8.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0789"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.09"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.73%  "
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.735.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.361.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.43%  "
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.70"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "45.795.80"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.65%  "
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.64%  "
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "243.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "38.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -11.44%  "
$ws.Range("E28").Value = "  -4.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.70"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "21.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.77"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +18.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.74"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.95%  "
$ws.Range("E33").Value = "  -4.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "147.40"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0770"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.31%  "
$ws.Range("E36").Value = "  -0.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.92"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.98%  "
$ws.Range("E38").Value = "  -2.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.08"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.94%  "
$ws.Range("E40").Value = "  -6.27%  "
$ws.Range("E41").Value = "  -2.43%  "
$ws.Range("E42").Value = "  -7.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.946.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "95.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -11.87%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.48"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.46%  "
$ws.Range("E48").Value = "  -6.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "98.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.607.64"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "68.71"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.42%  "
